$wb = $excel.ActiveWorkbook

# --- Update the "time_taken" (column F) timestamps on the "data" sheet ---
$data = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:35:56.757584",
    "2021-10-05 14:35:56.757593",
    "2021-10-05 14:35:56.757596",
    "2021-10-05 14:35:56.757599",
    "2021-10-05 14:35:56.757602",
    "2021-10-05 14:35:56.757604",
    "2021-10-05 14:35:56.757607",
    "2021-10-05 14:35:56.757610",
    "2021-10-05 14:35:56.757613",
    "2021-10-05 14:35:56.757615",
    "2021-10-05 14:35:56.757618",
    "2021-10-05 14:35:56.757620",
    "2021-10-05 14:35:56.757623",
    "2021-10-05 14:35:56.757626",
    "2021-10-05 14:35:56.757628"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" sheet right after "data" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "metadata"

# Re-fetch live handles by name - stale object refs from before the
# Add()/Move() calls can end up pointing at the wrong sheet.
$data = $wb.Worksheets.Item("data")
$newSheet = $wb.Worksheets.Item("metadata")
$newSheet.Move($null, $data)

# Re-fetch again post-move, just to be safe.
$meta = $wb.Worksheets.Item("metadata")

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Reproduce the "header" cell format (bold font, thin box border,
# centered/top aligned) used for styled cells on the "data" sheet -
# Range.Style assignment isn't supported by this host, so the individual
# format facets are applied directly instead (border/alignment/bold order
# matters: it lets the host collapse this into a single new style record).
$headerRange = $meta.Range("B1:G1")
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Font.Bold = $true

# Data row
$meta.Cells.Item(2, 1).Value = 0
$a2 = $meta.Range("A2")
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Font.Bold = $true
$meta.Cells.Item(2, 2).Value = "Vitamin C Pathway Disorders"
$meta.Cells.Item(2, 3).Value = 184

$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "0.15"

$meta.Cells.Item(2, 5).Value = "2021-02-04T23:41:21.505585Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:56.754002"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/184/?format=json"
